# Append 24 new log rows (539-562) to the worksheet, matching the
# "many small improvements" data refresh captured in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(539, "2022-03-22 07:31:14", 3600, 1647930674, 58, 885, 31)
    ,@(540, "2022-03-22 07:52:18", 3600, 1647931939, 60, 680, 40)
    ,@(541, "2022-03-22 07:52:48", 3600, 1647931939, 60, 680, 40)
    ,@(542, "2022-03-22 07:53:18", 3600, 1647931939, 60, 680, 40)
    ,@(543, "2022-03-22 07:53:58", 3600, 1647932039, 60, -721383, 43)
    ,@(544, "2022-03-22 08:10:21", 3600, 1647933022, 58, 186, 46)
    ,@(545, "2022-03-22 08:21:22", 3600, 1647933683, 58, 187, 52)
    ,@(546, "2022-03-22 08:28:24", 3600, 1647934105, 60, 178, 55)
    ,@(547, "2022-03-22 09:23:31", 3600, 1647937412, 52, 242, 68)
    ,@(548, "2022-03-22 09:52:35", 3600, 1647939156, 54, 184, 71)
    ,@(549, "2022-03-22 09:57:36", 3600, 1647939457, 60, 186, 74)
    ,@(550, "2022-03-22 10:03:37", 3600, 1647939818, 58, 185, 77)
    ,@(551, "2022-03-22 10:26:39", 3600, 1647941200, 60, 191, 81)
    ,@(552, "2022-03-22 10:29:40", 3600, 1647941381, 60, 187, 84)
    ,@(553, "2022-03-22 10:43:43", 3600, 1647942224, 60, 185, 90)
    ,@(554, "2022-03-22 11:01:45", 3600, 1647943307, 62, 179, 96)
    ,@(555, "2022-03-22 11:06:47", 3600, 1647943608, 56, 188, 99)
    ,@(556, "2022-03-22 11:17:48", 3600, 1647944270, 62, 200, 102)
    ,@(557, "2022-03-22 11:20:10", 3600, 1647944412, 62, 189, 102)
    ,@(558, "2022-03-22 11:23:49", 3600, 1647944631, 58, 198, 105)
    ,@(559, "2022-03-22 11:54:53", 3600, 1647946494, 62, 192, 108)
    ,@(560, "2022-03-22 11:59:53", 3600, 1647946795, 60, 203, 111)
    ,@(561, "2022-03-22 12:40:58", 3600, 1647949259, 58, 184, 114)
    ,@(562, "2022-03-22 12:43:57", 3600, 1647949439, 58, 222, 117)
)

$category = "8: 255`n"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $category
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = "10F872226797"
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    # Newline inside column B's text makes Excel auto-grow the row
    # height; restore the default row height/flag via AutoFit so the
    # row matches the rest of the sheet (no custom height stored).
    $ws.Rows.Item($r).AutoFit()
}
